$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed data cells (values scraped/refreshed in this automatic run)
$ws.Range('E2').Value = '2026-03-01 04:18:43'
$ws.Range('N2').Value = '-1.4 °C 3:58 TU'
$ws.Range('O2').Value = '-0.6 °C'
$ws.Range('E3').Value = '2026-03-01 04:18:46'
$ws.Range('E4').Value = '2026-03-01 04:18:49'
$ws.Range('H4').Value = "'97%"
$ws.Range('O4').Value = '8.4 °C'
$ws.Range('E5').Value = '2026-03-01 04:18:51'
$ws.Range('O5').Value = '-3.5 °C'
$ws.Range('E6').Value = '2026-03-01 04:18:54'
$ws.Range('H6').Value = "'85%"
$ws.Range('N6').Value = '8.8 °C 3:57 TU'
$ws.Range('O6').Value = '9.4 °C'
$ws.Range('E7').Value = '2026-03-01 04:18:57'
$ws.Range('N7').Value = '13.0 °C 3:52 TU'
$ws.Range('O7').Value = '13.3 °C'
$ws.Range('E8').Value = '2026-03-01 04:19:00'
$ws.Range('J8').Value = '1025.5 hPa'
$ws.Range('N8').Value = '9.3 °C 3:58 TU'
$ws.Range('E9').Value = '2026-03-01 04:19:02'
$ws.Range('H9').Value = "'61%"
$ws.Range('O9').Value = '11.7 °C'
$ws.Range('E10').Value = '2026-03-01 04:19:05'
$ws.Range('H10').Value = "'97%"
$ws.Range('N10').Value = '6.0 °C 3:59 TU'
$ws.Range('O10').Value = '6.8 °C'
$ws.Range('E11').Value = '2026-03-01 04:19:08'
$ws.Range('N11').Value = '6.1 °C 3:59 TU'
$ws.Range('E12').Value = '2026-03-01 04:19:10'
$ws.Range('H12').Value = "'72%"
$ws.Range('N12').Value = '9.3 °C 3:53 TU'
$ws.Range('O12').Value = '10.6 °C'
$ws.Range('E13').Value = '2026-03-01 04:19:13'
$ws.Range('N13').Value = '4.1 °C 3:31 TU'
$ws.Range('O13').Value = '4.4 °C'
$ws.Range('E14').Value = '2026-03-01 04:19:16'
$ws.Range('L14').Value = '9.7 km/h - 299º 3:57 TU'
$ws.Range('N14').Value = '10.3 °C 3:30 TU'
$ws.Range('O14').Value = '11.1 °C'
$ws.Range('E15').Value = '2026-03-01 04:19:18'
$ws.Range('O15').Value = '8.3 °C'
$ws.Range('E16').Value = '2026-03-01 04:19:21'
$ws.Range('N16').Value = '-5.7 °C 3:54 TU'
$ws.Range('O16').Value = '-4.7 °C'
$ws.Range('E17').Value = '2026-03-01 04:19:23'
$ws.Range('L17').Value = '8.6 km/h - 231º 3:48 TU'
$ws.Range('N17').Value = '1.1 °C 3:45 TU'
$ws.Range('O17').Value = '1.4 °C'
$ws.Range('E18').Value = '2026-03-01 04:19:25'
$ws.Range('H18').Value = "'99%"
$ws.Range('O18').Value = '7.1 °C'
$ws.Range('E19').Value = '2026-03-01 04:19:28'
$ws.Range('N19').Value = '5.9 °C 3:30 TU'
$ws.Range('E20').Value = '2026-03-01 04:19:30'
$ws.Range('N20').Value = '-3.6 °C 3:51 TU'
$ws.Range('O20').Value = '-2.9 °C'
$ws.Range('E21').Value = '2026-03-01 04:19:33'
$ws.Range('J21').Value = '1025.3 hPa'
$ws.Range('N21').Value = '6.2 °C 3:49 TU'
$ws.Range('E22').Value = '2026-03-01 04:19:36'
$ws.Range('O22').Value = '-5.1 °C'
$ws.Range('E23').Value = '2026-03-01 04:19:38'
$ws.Range('N23').Value = '-3.9 °C 3:58 TU'
$ws.Range('O23').Value = '-3.4 °C'
$ws.Range('E24').Value = '2026-03-01 04:19:41'
$ws.Range('E25').Value = '2026-03-01 04:19:44'
$ws.Range('H25').Value = "'96%"
$ws.Range('N25').Value = '-2.6 °C 3:43 TU'
$ws.Range('O25').Value = '-2.0 °C'
$ws.Range('E26').Value = '2026-03-01 04:19:46'
$ws.Range('N26').Value = '2.5 °C 3:31 TU'
$ws.Range('E27').Value = '2026-03-01 04:19:49'
$ws.Range('E28').Value = '2026-03-01 04:19:52'
$ws.Range('J28').Value = '1025.5 hPa'
$ws.Range('E29').Value = '2026-03-01 04:19:55'
$ws.Range('N29').Value = '8.9 °C 3:59 TU'
$ws.Range('O29').Value = '9.5 °C'
$ws.Range('E30').Value = '2026-03-01 04:19:57'
$ws.Range('H30').Value = "'78%"
$ws.Range('M30').Value = '11.6 °C 3:36 TU'
$ws.Range('O30').Value = '10.2 °C'
$ws.Range('E31').Value = '2026-03-01 04:20:00'
$ws.Range('N31').Value = '10.7 °C 3:39 TU'
$ws.Range('E32').Value = '2026-03-01 04:20:02'
$ws.Range('M32').Value = '3.7 °C 3:59 TU'
$ws.Range('O32').Value = '1.9 °C'
$ws.Range('E33').Value = '2026-03-01 04:20:05'
$ws.Range('E34').Value = '2026-03-01 04:20:07'
$ws.Range('N34').Value = '-0.3 °C 3:36 TU'
$ws.Range('E35').Value = '2026-03-01 04:20:10'
$ws.Range('E36').Value = '2026-03-01 04:20:12'
$ws.Range('H36').Value = "'78%"
$ws.Range('J36').Value = '1025.3 hPa'
$ws.Range('M36').Value = '11.7 °C 3:41 TU'
$ws.Range('O36').Value = '9.5 °C'
$ws.Range('E37').Value = '2026-03-01 04:20:15'
$ws.Range('N37').Value = '6.1 °C 3:50 TU'
$ws.Range('E38').Value = '2026-03-01 04:20:18'
$ws.Range('O38').Value = '8.8 °C'
$ws.Range('E39').Value = '2026-03-01 04:20:20'
$ws.Range('L39').Value = '27.7 km/h - 191º 3:47 TU'
$ws.Range('N39').Value = '-3.3 °C 3:58 TU'
$ws.Range('O39').Value = '-2.8 °C'
$ws.Range('E40').Value = '2026-03-01 04:20:23'
$ws.Range('G40').Value = '4 cm'
$ws.Range('H40').Value = "'87%"
$ws.Range('N40').Value = '6.8 °C 3:45 TU'
$ws.Range('O40').Value = '7.3 °C'
$ws.Range('E41').Value = '2026-03-01 04:20:26'
$ws.Range('N41').Value = '11.6 °C 3:58 TU'
$ws.Range('E42').Value = '2026-03-01 04:20:29'
$ws.Range('H42').Value = "'83%"
$ws.Range('N42').Value = '7.1 °C 3:59 TU'
$ws.Range('O42').Value = '9.2 °C'
$ws.Range('E43').Value = '2026-03-01 04:20:32'
$ws.Range('L43').Value = '5.0 km/h - 250º 3:49 TU'
$ws.Range('E44').Value = '2026-03-01 04:20:34'
$ws.Range('N44').Value = '-3.0 °C 3:54 TU'
$ws.Range('O44').Value = '-2.5 °C'
$ws.Range('E45').Value = '2026-03-01 04:20:37'
$ws.Range('J45').Value = '1027.0 hPa'
$ws.Range('E46').Value = '2026-03-01 04:20:40'
$ws.Range('O46').Value = '7.6 °C'

# Keep originally-empty cells empty (no data available for these fields)
$ws.Range('J2').Value = $null
$ws.Range('L2').Value = $null
$ws.Range('J3').Value = $null
$ws.Range('G4').Value = $null
$ws.Range('J5').Value = $null
$ws.Range('G6').Value = $null
$ws.Range('G7').Value = $null
$ws.Range('G8').Value = $null
$ws.Range('G9').Value = $null
$ws.Range('J9').Value = $null
$ws.Range('G10').Value = $null
$ws.Range('J10').Value = $null
$ws.Range('G11').Value = $null
$ws.Range('J11').Value = $null
$ws.Range('K11').Value = $null
$ws.Range('L11').Value = $null
$ws.Range('G12').Value = $null
$ws.Range('J12').Value = $null
$ws.Range('K12').Value = $null
$ws.Range('L12').Value = $null
$ws.Range('G14').Value = $null
$ws.Range('J14').Value = $null
$ws.Range('G15').Value = $null
$ws.Range('J15').Value = $null
$ws.Range('K15').Value = $null
$ws.Range('L15').Value = $null
$ws.Range('J16').Value = $null
$ws.Range('J17').Value = $null
$ws.Range('G18').Value = $null
$ws.Range('G19').Value = $null
$ws.Range('J19').Value = $null
$ws.Range('J20').Value = $null
$ws.Range('G21').Value = $null
$ws.Range('J22').Value = $null
$ws.Range('J23').Value = $null
$ws.Range('G24').Value = $null
$ws.Range('J25').Value = $null
$ws.Range('J27').Value = $null
$ws.Range('G28').Value = $null
$ws.Range('G29').Value = $null
$ws.Range('J29').Value = $null
$ws.Range('G30').Value = $null
$ws.Range('G31').Value = $null
$ws.Range('J32').Value = $null
$ws.Range('G33').Value = $null
$ws.Range('J34').Value = $null
$ws.Range('G36').Value = $null
$ws.Range('G37').Value = $null
$ws.Range('K37').Value = $null
$ws.Range('G38').Value = $null
$ws.Range('J38').Value = $null
$ws.Range('J39').Value = $null
$ws.Range('K40').Value = $null
$ws.Range('L40').Value = $null
$ws.Range('G41').Value = $null
$ws.Range('G42').Value = $null
$ws.Range('J42').Value = $null
$ws.Range('K42').Value = $null
$ws.Range('L42').Value = $null
$ws.Range('G43').Value = $null
$ws.Range('J43').Value = $null
$ws.Range('J44').Value = $null
$ws.Range('G46').Value = $null
